# Code Instructions.docx - "Add files via upload" edit
#
# Original single run:
#   "in other to import that correct data given the analyzed database,
#    the correct file name should be filled in line 12."
#
# becomes six runs (typo fixes "other"->"order", "that"->"the") plus a
# brand-new trailing sentence about the extra genomic-strand file:
#   "in " | "order" | " to import " | "the" |
#   " correct data given the analyzed database, the correct file name
#    should be filled in line 12." |
#   " Note that the extra file that are needed for the genomic strand
#    data are too large for Github so they are only available upon
#    request."

$d = $word.ActiveDocument

# Locate the sentence that needs to be rewritten/extended.
$target = $d.Content
$found = $target.Find.Execute(
    "in other to import that correct data given the analyzed database, the correct file name should be filled in line 12.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the target sentence to edit."
}

$start = $target.Start

# Wipe the old run's text (its empty <w:r> is dropped on save) and
# rebuild the sentence word-by-word/phrase-by-phrase so each piece lands
# in its own <w:r>, exactly like a human retyping parts of the sentence
# in Word would. Each InsertAfter() is issued on a freshly-collapsed
# Range anchored at the end of the previous insertion, which keeps the
# new runs distinct instead of being coalesced back into one.
$target.Text = ""

$r = $d.Range($start, $start)
$r.InsertAfter("in ")

$r = $d.Range($r.End, $r.End)
$r.InsertAfter("order")

$r = $d.Range($r.End, $r.End)
$r.InsertAfter(" to import ")

$r = $d.Range($r.End, $r.End)
$r.InsertAfter("the")

$r = $d.Range($r.End, $r.End)
$r.InsertAfter(" correct data given the analyzed database, the correct file name should be filled in line 12.")

$r = $d.Range($r.End, $r.End)
$r.InsertAfter(" Note that the extra file that are needed for the genomic strand data are too large for Github so they are only available upon request.")

Write-Output $d.Content.Text.Substring($start - 10, 330)
